# Updating working, added treatments, code cleanup.
# Insert a new "calendarWord" treatment group (treatment_id 4) into the
# "Web Parameters" sheet, mirroring the existing "calendarBar" group
# (rows 8-10), and bump up the treatment_id of all the later rows by one
# since a new group was inserted ahead of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")

# Insert three new blank rows before row 11 -- everything that used to be
# row 11 onward shifts down to row 14 onward.
$ws.Rows("11:13").Insert()

# New comment string used by the three new rows.
$newComment = "Calendar MEL question with word and no interaction month view."

# Row 11 (new): treatment_id 4, position 1 -- mirrors old row 8.
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = "calendarWord"
$ws.Cells.Item(11, 4).Value = "none"
$ws.Cells.Item(11, 5).Value = "none"
$ws.Cells.Item(11, 6).Value = 300
$ws.Cells.Item(11, 8).Value = 44593
$ws.Cells.Item(11, 9).Value = 700
$ws.Cells.Item(11, 11).Value = 44614
$ws.Cells.Item(11, 12).Value = 1100
$ws.Cells.Item(11, 14).Value = 100
$ws.Cells.Item(11, 15).Value = 100
$ws.Cells.Item(11, 20).Value = 8
$ws.Cells.Item(11, 21).Value = 8
$ws.Cells.Item(11, 22).Value = $newComment

# Row 12 (new): treatment_id 4, position 2 -- mirrors old row 9.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = "calendarWord"
$ws.Cells.Item(12, 4).Value = "none"
$ws.Cells.Item(12, 5).Value = "none"
$ws.Cells.Item(12, 6).Value = 500
$ws.Cells.Item(12, 8).Value = 44621
$ws.Cells.Item(12, 9).Value = 800
$ws.Cells.Item(12, 11).Value = 44632
$ws.Cells.Item(12, 12).Value = 1100
$ws.Cells.Item(12, 14).Value = 100
$ws.Cells.Item(12, 15).Value = 100
$ws.Cells.Item(12, 20).Value = 8
$ws.Cells.Item(12, 21).Value = 8
$ws.Cells.Item(12, 22).Value = $newComment

# Row 13 (new): treatment_id 4, position 3 -- mirrors old row 10.
$ws.Cells.Item(13, 1).Value = 4
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(13, 3).Value = "calendarWord"
$ws.Cells.Item(13, 4).Value = "none"
$ws.Cells.Item(13, 5).Value = "none"
$ws.Cells.Item(13, 6).Value = 300
$ws.Cells.Item(13, 8).Value = 44652
$ws.Cells.Item(13, 9).Value = 1000
$ws.Cells.Item(13, 11).Value = 44666
$ws.Cells.Item(13, 12).Value = 1100
$ws.Cells.Item(13, 14).Value = 100
$ws.Cells.Item(13, 15).Value = 100
$ws.Cells.Item(13, 20).Value = 8
$ws.Cells.Item(13, 21).Value = 8
$ws.Cells.Item(13, 22).Value = $newComment

# All of the rows that used to start at row 11 (now at row 14 onward, down
# through the last data row 21) get their treatment_id (column A)
# incremented by one, since a new treatment group was inserted before them.
for ($r = 14; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# Update the sheet view: scroll back to the left and select B14.
$ws.Activate()
$ws.Range("B14").Select()
